$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
Write-Host $ws.Name
